$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: rename the Windows demo VM / CVE row ----------------------
# Copy the existing row-2 cell formatting across to row 3 first (range
# copy, not whole-row copy, so it only touches the 5 used columns -
# column F on row 2 is an unused/empty styled cell that row 3 doesn't
# have in the target layout).
$ws.Range("A2:E2").Copy()
$ws.Range("A3:E3").PasteSpecial(-4122)

# New, shorter PowerShell command text for the Windows/Python row (the
# trailing "Optional: Verify installation" / "python --version" lines
# were dropped).
$windowsCommands = @'
# Define version and download path
$pythonVersion = "3.11.9"
$installerUrl = "https://www.python.org/ftp/python/$pythonVersion/python-$pythonVersion-amd64.exe"
$installerPath = "$env:TEMP\python-$pythonVersion-amd64.exe"
 
# Download the installer
Invoke-WebRequest -Uri $installerUrl -OutFile $installerPath -UseBasicParsing
 
# Run the installer silently
Start-Process -FilePath $installerPath -ArgumentList "/quiet InstallAllUsers=1 PrependPath=1 Include_test=0" -Wait
'@

$ws.Range("B2").Value = "windows-demo"
$ws.Range("C2").Value = $windowsCommands
$ws.Range("D2").Value = "CVE-Shell-Test-Windows"
$ws.Range("E2").Value = "update python version to 3.11.9"

$ws.Rows.Item(2).RowHeight = 210

# --- Row 3: brand-new Linux/openssl demo row ---------------------------
$linuxCommands = @'
# Update package lists
sudo apt update
# Install a specific version of openssl (3.0.14)
sudo apt install openssl=3.0.14-1~$(lsb_release -cs)1
# Hold the package at this version (optional)
sudo apt-mark hold openssl
'@

$ws.Range("A3").Value = "AZ-AS-SUB-EX-N-SEQ02125-CORE"
$ws.Range("B3").Value = "linux-demo"
$ws.Range("C3").Value = $linuxCommands
$ws.Range("D3").Value = "CVE-Shell-Test-Ubuntu"
$ws.Range("E3").Value = "Update openssl version to 3.0.14"

$ws.Rows.Item(3).RowHeight = 120

# --- Selection moves from D3 to E3 -------------------------------------
$ws.Range("E3").Select()
